# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# F2: 1577 -> 1580
# F3: 100  -> 104
# F4: 35   -> 36

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1580
    $ws.Range("F3").Value = 104
    $ws.Range("F4").Value = 36
}
